$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.983.48"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "1.742.33"
$ws.Range("E3").Value = "  -0.41%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.50"
$ws.Range("E5").Value = "  +7.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5157"
$ws.Range("E7").Value = "  -2.57%  "

$ws.Range("E8").Value = "  -0.77%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06195"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "1.742.39"
$ws.Range("E10").Value = "  -0.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07233"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("E12").Value = "  -1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6525"
$ws.Range("E13").Value = "  +1.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.638"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.77"
$ws.Range("E15").Value = "  -0.90%  "

$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("D18").Value = "26.001.27"
$ws.Range("E18").Value = "  +0.09%  "

$ws.Range("E19").Value = "  +1.96%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006817"
$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").Value = "1.964.99"
$ws.Range("E21").Value = "  -0.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.283"
$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.691"
$ws.Range("E23").Value = "  -1.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.378"
$ws.Range("E24").Value = "  +2.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.29"
$ws.Range("E25").Value = "  -2.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.512"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.781"
$ws.Range("E28").Value = "  -1.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.98"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.972"
$ws.Range("E30").Value = "  +4.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08226"
$ws.Range("E31").Value = "  -0.98%  "

$ws.Range("E32").Value = "  -0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04705"
$ws.Range("E33").Value = "  +3.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.656"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6254"
$ws.Range("E36").Value = "  -1.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.733"
$ws.Range("E37").Value = "  +0.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01617"
$ws.Range("E38").Value = "  +1.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.923"
$ws.Range("E39").Value = "  -0.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9997"
$ws.Range("E40").Value = "  +0.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.64"
$ws.Range("E41").Value = "  +2.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7631"
$ws.Range("E42").Value = "  +3.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3865"
$ws.Range("E43").Value = "  -1.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.036"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.342"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("E46").Value = "  -0.48%  "

$ws.Range("E47").Value = "  +2.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05230"
$ws.Range("E48").Value = "  -2.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.86"
$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.566"
$ws.Range("E50").Value = "  -1.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3437"
$ws.Range("E51").Value = "  -0.85%  "
